$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column layout: A=Datum, B=Pocet.potvrdenych.PCR.testami, C=Dennych.PCR.testov,
# D=Dennych.PCR.prirastkov, E=Pocet.umrti, F=AgTests, G=AgPosit

# Update existing rows (column F = AgTests, column G = AgPosit)
$ws.Range("F426").Value = 107041
$ws.Range("G426").Value = 381

$ws.Range("F427").Value = 90454

$ws.Range("F428").Value = 102372

$ws.Range("F429").Value = 178092

$ws.Range("F431").Value = 171514

$ws.Range("F432").Value = 123375

$ws.Range("F435").Value = 82935

$ws.Range("F436").Value = 144871

$ws.Range("F440").Value = 73556

$ws.Range("F445").Value = 84586

$ws.Range("F446").Value = 86183

$ws.Range("F447").Value = 67318

$ws.Range("F448").Value = 61506

$ws.Range("F449").Value = 59823

$ws.Range("F451").Value = 85647

$ws.Range("F452").Value = 74607

$ws.Range("F453").Value = 70078

$ws.Range("F454").Value = 51671

$ws.Range("F455").Value = 50443

$ws.Range("F457").Value = 77187
$ws.Range("G457").Value = 131

$ws.Range("F458").Value = 67810
$ws.Range("G458").Value = 73

$ws.Range("F459").Value = 57081
$ws.Range("G459").Value = 83

# Append new row 460 with the latest daily stats entry
$ws.Range("A460").NumberFormat = "yyyy-mm-dd"
$ws.Range("A460").Value = 44354
$ws.Range("B460").Value = 390546
$ws.Range("C460").Value = 5776
$ws.Range("D460").Value = 95
$ws.Range("E460").Value = 12423
$ws.Range("F460").Value = 47165
$ws.Range("G460").Value = 150

$wb.Save()
